$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.470.78"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "2.249.54"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.66"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.98"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.573"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.89"
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0812"
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.19"
$ws.Range("E12").Value = "  -2.05%  "
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").Value = "2.357.18"
$ws.Range("E14").Value = "  +4.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.838"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.65"
$ws.Range("E16").Value = "  -1.67%  "
$ws.Range("D17").Value = "44.196.67"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.42"
$ws.Range("E19").Value = "  -4.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.38"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "65.82"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.32"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.96"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "39.02"
$ws.Range("E26").Value = "  +6.31%  "
$ws.Range("E27").Value = "  +3.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.86"
$ws.Range("E28").Value = "  -3.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.98"
$ws.Range("E29").Value = "  -4.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.07"
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.44"
$ws.Range("E31").Value = "  -2.69%  "
$ws.Range("E32").Value = "  -3.52%  "
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.14"
$ws.Range("E34").Value = "  -10.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.110"
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("E37").Value = "  -3.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.49"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.84"
$ws.Range("E39").Value = "  -7.47%  "
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0303"
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("D43").Value = "1.738.23"
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "81.51"
$ws.Range("E44").Value = "  -6.86%  "
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.63"
$ws.Range("E46").Value = "  +4.94%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.04"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.95"
$ws.Range("E48").Value = "  -4.26%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "14.78"
$ws.Range("E49").Value = "  +5.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.19"
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.50"
$ws.Range("E51").Value = "  -0.11%  "
